$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 3.123541145015474
$ws.Cells.Item(2, 5).Value = 2.651224325850143
$ws.Cells.Item(3, 3).Value = 3.959010658874851
$ws.Cells.Item(3, 5).Value = 2.48038747547068
$ws.Cells.Item(4, 3).Value = 4.722695063536686
$ws.Cells.Item(4, 5).Value = 3.580489688170352
$ws.Cells.Item(5, 3).Value = 6.739021039846627
$ws.Cells.Item(5, 5).Value = 2.457570659067509
$ws.Cells.Item(6, 3).Value = 2.619839412265601
$ws.Cells.Item(6, 5).Value = 1.929582042845213
$ws.Cells.Item(7, 3).Value = -0.7919564768266385
$ws.Cells.Item(7, 5).Value = 1.978648203842193
$ws.Cells.Item(8, 3).Value = 1.877689851450803
$ws.Cells.Item(8, 5).Value = 3.188468414048606
$ws.Cells.Item(9, 3).Value = 2.705004599189187
$ws.Cells.Item(9, 5).Value = 2.471779821159181
$ws.Cells.Item(10, 3).Value = 1.110374544249249
$ws.Cells.Item(10, 5).Value = 2.958768964947134
$ws.Cells.Item(11, 3).Value = 2.267566233338814
$ws.Cells.Item(11, 5).Value = 2.824054578526636
$ws.Cells.Item(12, 3).Value = 2.688433258834588
$ws.Cells.Item(12, 5).Value = 2.558605038029849
$ws.Cells.Item(13, 3).Value = 1.014079695989589
$ws.Cells.Item(13, 5).Value = 2.1961756525541
$ws.Cells.Item(14, 3).Value = 3.013853578092252
$ws.Cells.Item(14, 5).Value = 2.721358666668894
$ws.Cells.Item(15, 3).Value = 1.331295149770684
$ws.Cells.Item(15, 5).Value = 1.104792173470215
$ws.Cells.Item(16, 3).Value = 0.04589006555719699
$ws.Cells.Item(16, 5).Value = 1.367746832546346
$ws.Cells.Item(17, 3).Value = 0.009546395482029624
$ws.Cells.Item(17, 5).Value = 0.9083591402091473
$ws.Cells.Item(18, 3).Value = 0.8709390141433015
$ws.Cells.Item(18, 5).Value = 1.177606443599988
$ws.Cells.Item(19, 3).Value = 0.7652063367885598
$ws.Cells.Item(19, 5).Value = 1.680477107215861
$ws.Cells.Item(20, 3).Value = 2.267579219134386
$ws.Cells.Item(20, 5).Value = 2.225111592343887
$ws.Cells.Item(21, 3).Value = 3.146753122914103
$ws.Cells.Item(21, 5).Value = 1.799793194919874
$ws.Cells.Item(22, 3).Value = 1.769033835366818
$ws.Cells.Item(22, 5).Value = 0.2856860139923256
$ws.Cells.Item(23, 3).Value = -4.774715709990263
$ws.Cells.Item(23, 5).Value = 0.754926127539246
$ws.Cells.Item(24, 3).Value = 1.95493704440024
$ws.Cells.Item(24, 5).Value = 2.928378677701393
$ws.Cells.Item(25, 3).Value = 3.478075069442799
$ws.Cells.Item(25, 5).Value = 1.768431385360159
$ws.Cells.Item(26, 3).Value = 1.232342134690434
$ws.Cells.Item(26, 5).Value = 1.134779475590464
$ws.Cells.Item(27, 3).Value = 0.2542811494408159
$ws.Cells.Item(27, 5).Value = 1.353526127153426
$ws.Cells.Item(28, 3).Value = 1.519778766382096
$ws.Cells.Item(28, 5).Value = 0.7212678493511149
$ws.Cells.Item(29, 3).Value = 1.469441753880329
$ws.Cells.Item(29, 5).Value = 1.350400980499855
$ws.Cells.Item(30, 3).Value = 1.638203081492495
$ws.Cells.Item(30, 5).Value = 1.323745783269614
$ws.Cells.Item(31, 3).Value = 2.268697431234346
$ws.Cells.Item(31, 5).Value = 2.455413743911294
$ws.Cells.Item(32, 3).Value = 1.984425467899631
$ws.Cells.Item(32, 5).Value = 0.6687400825358569
$ws.Cells.Item(33, 3).Value = 0.6066448776129052
$ws.Cells.Item(33, 5).Value = 0.6230021429014077
$ws.Cells.Item(34, 3).Value = -4.243076347305386
$ws.Cells.Item(34, 5).Value = -2.163103471150829
$ws.Cells.Item(35, 3).Value = 1.438499295329754
$ws.Cells.Item(35, 5).Value = 1.882119284761474
$ws.Cells.Item(36, 3).Value = 1.906593537051537
$ws.Cells.Item(36, 5).Value = 1.473394465200051
$ws.Cells.Item(37, 3).Value = 0.08348019664223827
$ws.Cells.Item(37, 5).Value = 0.8702074629614476
$ws.Cells.Item(38, 3).Value = -0.214505326882275
$ws.Cells.Item(38, 5).Value = 0.8092352694139215
$ws.Cells.Item(39, 3).Value = 0.1651547428133782
$ws.Cells.Item(39, 5).Value = 0.7536567386490001